# Fruta / hortaliza, semanal
# Insert two new weekly price records (2021-09-09, serial 44448) above the
# existing rows 170-174, pushing the previous data down to rows 172-176.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at position 170 (shifts old rows 170..174 to 172..176)
$ws.Rows.Item(170).Insert()
$ws.Rows.Item(170).Insert()

# New row 170: Kiwi Hayward - Especial
$ws.Cells.Item(170, 1).Value = 5
$ws.Cells.Item(170, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(170, 3).Value = "Maule"
$ws.Cells.Item(170, 4).Value = 44448
$ws.Cells.Item(170, 5).Value = 7
$ws.Cells.Item(170, 6).Value = "Fruta"
$ws.Cells.Item(170, 7).Value = 100101
$ws.Cells.Item(170, 8).Value = "Berries"
$ws.Cells.Item(170, 9).Value = 100101007
$ws.Cells.Item(170, 10).Value = "Kiwi"
$ws.Cells.Item(170, 11).Value = "Hayward"
$ws.Cells.Item(170, 12).Value = "Especial"
$ws.Cells.Item(170, 13).Value = 180
$ws.Cells.Item(170, 14).Value = 13000
$ws.Cells.Item(170, 15).Value = 13000
$ws.Cells.Item(170, 16).Value = 13000
$ws.Cells.Item(170, 17).Value = "`$/caja 18 kilos"
$ws.Cells.Item(170, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(170, 19).Value = 722
$ws.Cells.Item(170, 20).Value = 18

# New row 171: Kiwi Hayward - Primera
$ws.Cells.Item(171, 1).Value = 5
$ws.Cells.Item(171, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(171, 3).Value = "Maule"
$ws.Cells.Item(171, 4).Value = 44448
$ws.Cells.Item(171, 5).Value = 7
$ws.Cells.Item(171, 6).Value = "Fruta"
$ws.Cells.Item(171, 7).Value = 100101
$ws.Cells.Item(171, 8).Value = "Berries"
$ws.Cells.Item(171, 9).Value = 100101007
$ws.Cells.Item(171, 10).Value = "Kiwi"
$ws.Cells.Item(171, 11).Value = "Hayward"
$ws.Cells.Item(171, 12).Value = "Primera"
$ws.Cells.Item(171, 13).Value = 230
$ws.Cells.Item(171, 14).Value = 11000
$ws.Cells.Item(171, 15).Value = 11000
$ws.Cells.Item(171, 16).Value = 11000
$ws.Cells.Item(171, 17).Value = "`$/caja 18 kilos"
$ws.Cells.Item(171, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(171, 19).Value = 611
$ws.Cells.Item(171, 20).Value = 18
